$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Mean (S.D.) values in column C
$ws.Range("C2").Value = "0.785 (0.004)"
$ws.Range("C3").Value = "0.874 (0.005)"
$ws.Range("C4").Value = "0.878 (0.010)"
$ws.Range("C5").Value = "0.914 (0.011)"
$ws.Range("C6").Value = "0.939 (0.161)"
$ws.Range("C7").Value = "0.956 (0.045)"

# Swap CORnet S / CORnet Z labels in rows 8 and 9, and update their S.D. values
$ws.Range("B8").Value = "CORnet Z"
$ws.Range("C8").Value = "1.000 (0.000)"

$ws.Range("B9").Value = "CORnet S"
$ws.Range("C9").Value = "1.000 (0.001)"
